$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 27778834
$ws.Range("I19").Value = 802.125
$ws.Range("J19").Value = 50001260
$ws.Range("K19").Value = 802.125
$ws.Range("L19").Value = 50001260
$ws.Range("M19").Value = -627.125
$ws.Range("N19").Value = -50001610
$ws.Range("H33").Value = 33310342
$ws.Range("I33").Value = 49965290
$ws.Range("J33").Value = 450
$ws.Range("K33").Value = 49965290
$ws.Range("L33").Value = 450
$ws.Range("M33").Value = -49965061
$ws.Range("N33").Value = -908
$ws.Range("H40").Value = 1578.6666
$ws.Range("I40").Value = 2233.3333
$ws.Range("J40").Value = 1415
$ws.Range("K40").Value = 2233.3333
$ws.Range("L40").Value = 1415
$ws.Range("M40").Value = -2058.3333
$ws.Range("N40").Value = -1765
$ws.Range("H64").Value = 3995.6365
$ws.Range("J64").Value = 4269.6924
$ws.Range("L64").Value = 4269.6924
$ws.Range("N64").Value = -4765.6924
$ws.Range("H67").Value = 3995.6365
$ws.Range("J67").Value = 4269.6924
$ws.Range("L67").Value = 4269.6924
$ws.Range("N67").Value = -5985.6924
$ws.Range("H76").Value = 4106.522
$ws.Range("J76").Value = 5190
$ws.Range("L76").Value = 5190
$ws.Range("N76").Value = -5820
$ws.Range("H79").Value = 4106.522
$ws.Range("J79").Value = 5190
$ws.Range("L79").Value = 5190
$ws.Range("N79").Value = -7374
$ws.Range("H98").Value = 4882.5835
$ws.Range("I98").Value = 4950
$ws.Range("J98").Value = 4747.75
$ws.Range("K98").Value = 4950
$ws.Range("L98").Value = 4747.75
$ws.Range("M98").Value = -3452
$ws.Range("N98").Value = -7743.75
$ws.Range("H122").Value = 4882.5835
$ws.Range("I122").Value = 4950
$ws.Range("J122").Value = 4747.75
$ws.Range("K122").Value = 14850
$ws.Range("L122").Value = 14243.25
$ws.Range("M122").Value = -12400
$ws.Range("N122").Value = -19143.25
$ws.Range("H135").Value = 42857956
$ws.Range("I135").Value = 20833888
$ws.Range("J135").Value = 90910460
$ws.Range("K135").Value = 187504992
$ws.Range("L135").Value = 818194140
$ws.Range("M135").Value = -187502457
$ws.Range("N135").Value = -818199210
$ws.Range("H138").Value = 3908.9412
$ws.Range("I138").Value = 3733
$ws.Range("J138").Value = 3963.077
$ws.Range("K138").Value = 11199
$ws.Range("L138").Value = 11889.231
$ws.Range("M138").Value = -6059
$ws.Range("N138").Value = -22169.231

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 37090.5
$ws.Range("J62").Value = 37090.5
$ws.Range("L62").Value = 37090.5
$ws.Range("N62").Value = -38462.5
$ws.Range("H65").Value = 37090.5
$ws.Range("J65").Value = 37090.5
$ws.Range("L65").Value = 111271.5
$ws.Range("N65").Value = -118135.5
$ws.Range("H107").Value = 3109.3333
$ws.Range("I107").Value = 2968
$ws.Range("J107").Value = 3533.3333
$ws.Range("K107").Value = 2968
$ws.Range("L107").Value = 3533.3333
$ws.Range("M107").Value = -1048
$ws.Range("N107").Value = -7373.3333
$ws.Range("H134").Value = 65485.438
$ws.Range("I134").Value = 2671.0833
$ws.Range("K134").Value = 8013.249899999999
$ws.Range("M134").Value = -5478.249899999999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("H96").Value = 11500
$ws.Range("J96").Value = 11500
$ws.Range("L96").Value = 11500
$ws.Range("N96").Value = -16992
$ws.Range("H114").Value = 41000
$ws.Range("J114").Value = 41000
$ws.Range("L114").Value = 41000
$ws.Range("N114").Value = -49678
$ws.Range("H127").Value = 3806200
$ws.Range("J127").Value = 3806200
$ws.Range("L127").Value = 3806200
$ws.Range("N127").Value = -3816120
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("M30").ClearContents()
$ws.Range("N30").ClearContents()
$ws.Range("M128").ClearContents()
$ws.Range("N128").ClearContents()
$ws.Range("N130").ClearContents()

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 11118081
$ws.Range("I5").Value = 338.72726
$ws.Range("J5").Value = 41691870
$ws.Range("K5").Value = 1016.18178
$ws.Range("L5").Value = 125075610
$ws.Range("M5").Value = -904.18178
$ws.Range("N5").Value = -125075834
$ws.Range("H113").Value = 399.48
$ws.Range("I113").Value = 474.5778
$ws.Range("J113").Value = 338.03638
$ws.Range("K113").Value = 1423.7334
$ws.Range("L113").Value = 1014.10914
$ws.Range("M113").Value = 746.2665999999999
$ws.Range("N113").Value = -5354.10914
$ws.Range("H122").Value = 1280
$ws.Range("I122").Value = 280
$ws.Range("J122").Value = 1480
$ws.Range("K122").Value = 2520
$ws.Range("L122").Value = 13320
$ws.Range("M122").Value = -70
$ws.Range("N122").Value = -18220
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("H124").Value = 2213.8333
$ws.Range("J124").Value = 3664.3333
$ws.Range("L124").Value = 10992.9999
$ws.Range("N124").Value = -20812.9999
$ws.Range("H125").Value = 5257.5
$ws.Range("I125").Value = 3030
$ws.Range("K125").Value = 9090
$ws.Range("M125").Value = -4170
$ws.Range("H126").Value = 3486
$ws.Range("I126").Value = 1230
$ws.Range("J126").Value = 4050
$ws.Range("K126").Value = 3690
$ws.Range("L126").Value = 12150
$ws.Range("M126").Value = 1250
$ws.Range("N126").Value = -22030
$ws.Range("H129").Value = 1385.8387
$ws.Range("I129").Value = 760
$ws.Range("J129").Value = 1641.8636
$ws.Range("K129").Value = 2280
$ws.Range("L129").Value = 4925.5908
$ws.Range("M129").Value = 2720
$ws.Range("N129").Value = -14925.5908
$ws.Range("H131").Value = 25585.281
$ws.Range("I131").Value = 1068.75
$ws.Range("J131").Value = 31912.129
$ws.Range("K131").Value = 3206.25
$ws.Range("L131").Value = 95736.387
$ws.Range("M131").Value = 1833.75
$ws.Range("N131").Value = -105816.387
$ws.Range("H132").Value = 2949.3235
$ws.Range("I132").Value = 3250.5
$ws.Range("J132").Value = 2909.1667
$ws.Range("K132").Value = 29254.5
$ws.Range("L132").Value = 26182.5003
$ws.Range("M132").Value = -26724.5
$ws.Range("N132").Value = -31242.5003
$ws.Range("H133").Value = 4481.3
$ws.Range("I133").Value = 2299.0908
$ws.Range("J133").Value = 5744.684
$ws.Range("K133").Value = 6897.2724
$ws.Range("L133").Value = 17234.052
$ws.Range("M133").Value = -1837.2724
$ws.Range("N133").Value = -27354.052
$ws.Range("H134").Value = 5059.1924
$ws.Range("J134").Value = 3938.1765
$ws.Range("L134").Value = 11814.5295
$ws.Range("N134").Value = -21954.5295
$ws.Range("H135").Value = 11118081
$ws.Range("I135").Value = 338.72726
$ws.Range("J135").Value = 41691870
$ws.Range("K135").Value = 3048.54534
$ws.Range("L135").Value = 375226830
$ws.Range("M135").Value = -513.5453400000001
$ws.Range("N135").Value = -375231900
$ws.Range("H137").Value = 126874.5
$ws.Range("I137").Value = 2499.3333
$ws.Range("K137").Value = 7497.999899999999
$ws.Range("M137").Value = -2397.999899999999
$ws.Range("M123").ClearContents()

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 13363.363
$ws.Range("I122").Value = 15777.444
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 47332.33199999999
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -44882.33199999999
$ws.Range("N122").Value = -12400
$ws.Range("H128").Value = 61453.332
$ws.Range("J128").Value = 61453.332
$ws.Range("L128").Value = 61453.332
$ws.Range("N128").Value = -71413.33199999999
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I7").Value = 3000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2888
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6530
$ws.Range("H138").Value = 70000
$ws.Range("J138").Value = 70000
$ws.Range("L138").Value = 70000
$ws.Range("N138").Value = -80280
$ws.Range("H139").Value = 75797.5
$ws.Range("J139").Value = 75797.5
$ws.Range("L139").Value = 75797.5
$ws.Range("N139").Value = -86077.5
$ws.Range("N7").ClearContents()
$ws.Range("N126").ClearContents()

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 641.7368
$ws.Range("I113").Value = 381.375
$ws.Range("J113").Value = 711.1667
$ws.Range("K113").Value = 1144.125
$ws.Range("L113").Value = 2133.5001
$ws.Range("M113").Value = 1025.875
$ws.Range("N113").Value = -6473.5001
$ws.Range("H122").Value = 2258.3333
$ws.Range("I122").Value = 2009.091
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6027.272999999999
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -3577.272999999999
$ws.Range("N122").Value = -19900
$ws.Range("H133").Value = 52851
$ws.Range("J133").Value = 52851
$ws.Range("L133").Value = 52851
$ws.Range("N133").Value = -62971
$ws.Range("H136").Value = 4771.755
$ws.Range("I136").Value = 1406.9615
$ws.Range("K136").Value = 4220.8845
$ws.Range("M136").Value = -1670.8845

Write-Output "All edits applied."